# Applies the "Updated symbol list" price/volume refresh to cryptos.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values that are stored as text (e.g. "23.00", "0.001260")
# so a leading apostrophe is used to force Excel to keep them as text instead of
# converting them to numbers (which would also strip significant trailing zeros).
$ws.Range("D2").Value  = "'243.13"
$ws.Range("D3").Value  = "'23.04"
$ws.Range("D4").Value  = "'5.406"
$ws.Range("D5").Value  = "'0.05957"
$ws.Range("D6").Value  = "'3.427"
$ws.Range("D7").Value  = "'6.507"
$ws.Range("D8").Value  = "'0.8122"
$ws.Range("D9").Value  = "'0.9237"
$ws.Range("D10").Value = "'0.1431"
$ws.Range("D11").Value = "'0.07424"
$ws.Range("D12").Value = "'0.03287"
$ws.Range("D14").Value = "'0.09355"
$ws.Range("D16").Value = "'0.001584"
$ws.Range("D17").Value = "'0.04682"
$ws.Range("D18").Value = "'0.0005977"
$ws.Range("D19").Value = "'0.005852"

$ws.Range("D20").Value = "'0.001261"
$ws.Range("E20").Value = "19BitKanKANBestin24h"

$ws.Range("D21").Value = "'0.004795"
$ws.Range("D22").Value = "'0.00007992"
$ws.Range("D23").Value = "'3.574"
$ws.Range("D27").Value = "'0.0002339"
$ws.Range("D40").Value = "'0.03936"
$ws.Range("D41").Value = "'0.006354"

# Rows 42/43: CEJI and BKEXToken swap places (with refreshed prices/labels)
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1075"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002607"
$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").Value = "'0.008894"
$ws.Range("D45").Value = "'0.00005178"
$ws.Range("D47").Value = "'0.6798"
$ws.Range("D48").Value = "'0.002143"
